$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff. Plain decimal-looking values in
# column D must be forced to Text so Excel does not auto-convert them to
# floating point numbers (which would corrupt values like "603.46" into
# "603.46000000000004" and collapse "1.00" into "1"). We briefly mark the
# cell as Text, assign the literal string, then clear the format again so the
# cell keeps its original (default/general) style, matching the source file.

# Row 2
$ws.Range('D2').Value = '68.866.57'
$ws.Range('E2').Value = '  +0.98%  '
# Row 3
$ws.Range('D3').Value = '2.730.94'
$ws.Range('E3').Value = '  +3.44%  '
# Row 4
$ws.Range('E4').Value = '  +0.08%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.46'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.45%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.87'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +6.33%  '
# Row 8
$ws.Range('E8').Value = '  +0.89%  '
# Row 9
$ws.Range('D9').Value = '2.731.23'
$ws.Range('E9').Value = '  +3.49%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.145'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.57%  '
# Row 11
$ws.Range('E11').Value = '  +4.62%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.33'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.35%  '
# Row 13
$ws.Range('E13').Value = '  -0.19%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.63'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.98%  '
# Row 15
$ws.Range('E15').Value = '  +3.59%  '
# Row 16
$ws.Range('E16').Value = '  +1.71%  '
# Row 17
$ws.Range('D17').Value = '68.750.77'
$ws.Range('E17').Value = '  +1.01%  '
# Row 18
$ws.Range('D18').Value = '2.731.68'
$ws.Range('E18').Value = '  +3.40%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.86'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.61%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '372.65'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.61%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.66'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.68%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.53'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.36%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.95'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.78%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.52%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.61'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.90%  '
# Row 26
$ws.Range('E26').Value = '  +0.01%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.95'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.89%  '
# Row 28
$ws.Range('D28').Value = '2.872.01'
$ws.Range('E28').Value = '  +3.53%  '
# Row 29
$ws.Range('E29').Value = '  +2.75%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '588.06'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.66%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.07%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.29'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.35%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.44'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.94%  '
# Row 34
$ws.Range('E34').Value = '  +5.78%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.132'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.21%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.62'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.28%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '161.10'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.44%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.90'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.10%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.381'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.41%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.93'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.66%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.49'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.72%  '
# Row 43
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.99'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.14%  '
# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.66'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.36%  '
# Row 45
$ws.Range('E45').Value = '  +0.05%  '
# Row 46
$ws.Range('E46').Value = '  -2.22%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '41.09'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.38%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '156.58'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.02%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.96'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +6.19%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.79'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.95%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.603'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.05%  '
